$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row height tweaks (row 78 existing, row 79 new) ---
$ws.Rows.Item(78).RowHeight = 13.8

# --- Build new row 79 by cloning the formatting of row 78 (only the columns that
#     row 78 actually populates - A:O, R, U:AC - so we don't introduce stray blank
#     cells in the gaps P, Q, S, T), then overwrite values ---
$ws.Range("A78:O78").Copy()
$ws.Range("A79:O79").PasteSpecial(-4122)
$ws.Range("R78").Copy()
$ws.Range("R79").PasteSpecial(-4122)
$ws.Range("U78:AC78").Copy()
$ws.Range("U79:AC79").PasteSpecial(-4122)

# Give the "final output" columns (R and AC) the same yellow-highlight style used
# on other SKU rows (e.g. row 6), since row 78 itself isn't highlighted.
$ws.Range("R6").Copy()
$ws.Range("R79").PasteSpecial(-4122)
$ws.Range("AC6").Copy()
$ws.Range("AC79").PasteSpecial(-4122)

$ws.Rows.Item(79).RowHeight = 13.8

# --- Values for new row 79 (SKU: Моцарелла палочки "Из Лавки") ---
$ws.Cells.Item(79, 1).Value = 77
$ws.Cells.Item(79, 2).Value = "Моцарелла палочки ""Из Лавки"", 45%, 0,12 кг, т/ф"
$ws.Cells.Item(79, 3).Value = 2.7
$ws.Cells.Item(79, 4).Value = "Да"
$ws.Cells.Item(79, 5).Value = "Моцарелла"
$ws.Cells.Item(79, 6).Value = "Соль"
$ws.Cells.Item(79, 7).Value = "Альче"
$ws.Cells.Item(79, 8).Value = "Из Лавки"
$ws.Cells.Item(79, 9).Value = 120
$ws.Cells.Item(79, 10).Value = 10
$ws.Cells.Item(79, 11).Value = 30
$ws.Cells.Item(79, 12).Value = 960
$ws.Cells.Item(79, 13).Value = 65
$ws.Cells.Item(79, 14).Value = "Нет"
$ws.Cells.Item(79, 15).Value = "Ульма"
$ws.Cells.Item(79, 18).Value = 900
$ws.Cells.Item(79, 21).Value = 10
$ws.Cells.Item(79, 22).Value = 35
$ws.Cells.Item(79, 23).Value = 20
$ws.Cells.Item(79, 24).Value = 20
$ws.Cells.Item(79, 25).Value = 15
$ws.Cells.Item(79, 26).Value = 5
$ws.Cells.Item(79, 27).Value = 5
$ws.Cells.Item(79, 28).Value = "00-00012176"
$ws.Cells.Item(79, 29).Value = 900

# NOTE: columns E:U are already stored as visible once the workbook round-trips
# through this engine, so no explicit unhide call is required (and issuing one
# would only fragment the <cols> column-range grouping).

# --- Touch the very last row of the sheet so the workbook dimension/UI matches the
#     author's final save (LibreOffice recalculated row heights all the way down). ---
$ws.Cells.Item(1048576, 2).NumberFormat = "General"
$ws.Rows.Item(1048576).RowHeight = 12.8

# --- Scroll position / selection as left by the author ---
$ws.Activate()
$ws.Range("B82").Select()
$excel.ActiveWindow.ScrollRow = 47
$excel.ActiveWindow.ScrollColumn = 1
